# Daily attendance processing - 2026-01-12 03:40:22
# Swap the order of names in the "Recorded By" column (G) from
# "System, dnasr281@gmail.com" to "dnasr281@gmail.com, System"
# for every row where that exact text is present.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetRows = @(2,3,4,5,6,7,8,16,17,22,23,24,26,29,37,38,43,44,45,47,50,58,59,64,65,66,68,71,79,80,85,86,87,88,89,90,91,99,100,105,106,107,108,109,110,111,119,120,125,126,127,128,129,130,131,139,140,145,146,147,148,149,150,151,159,160,165,166,167,168,169,170,171,179,180,185,186,187,189,192,200,201,206,207,208,210,213,221,222,227,228,229,231,234,242,243)

foreach ($row in $targetRows) {
    $cell = $ws.Cells.Item($row, 7)
    if ($cell.Value() -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
}
